$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 115677.51
$ws.Range("D2").Value = 258

$ws.Range("B3").Value = 953108.38
$ws.Range("C3").Value = 723.9357676353857
$ws.Range("D3").Value = 2141

$ws.Range("B4").Value = 1786705.19
$ws.Range("C4").Value = 87.46086253065994
$ws.Range("D4").Value = 2716

$ws.Range("B5").Value = 2885974.02
$ws.Range("C5").Value = 61.52491391151107
$ws.Range("D5").Value = 3342

$ws.Range("B6").Value = 4517432.77
$ws.Range("C6").Value = 56.53061111062945
$ws.Range("D6").Value = 4662

$ws.Range("B7").Value = 1774017.75
$ws.Range("C7").Value = -60.72951518435104
$ws.Range("D7").Value = 1680
